# Femacal de La Calera - Achicoria: insert a new weekly record.
# A new row of data (dated 44498) is inserted right before the existing
# row 41 (dated 44218), pushing all subsequent rows (old 41..136) down by
# one position (new 42..137), so the sheet grows from 136 to 137 data rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 41; Excel shifts rows 41-136 down to
# 42-137 and copies formatting (incl. the date number format) from the
# row above, just like a normal Excel "Insert Row" operation.
$ws.Rows(41).Insert()

# Populate the newly inserted row with the new record's data.
$ws.Range("A41").Value = 3
$ws.Range("B41").Value = "Femacal de La Calera"
$ws.Range("C41").Value = "Coquimbo"
$ws.Range("D41").Value = 44498
$ws.Range("E41").Value = 5
$ws.Range("F41").Value = 100112010
$ws.Range("G41").Value = "Achicoria"
$ws.Range("H41").Value = "Sin especificar"
$ws.Range("I41").Value = "Primera"
$ws.Range("J41").Value = 115
$ws.Range("K41").Value = 6000
$ws.Range("L41").Value = 6300
$ws.Range("M41").Value = 6143
$ws.Range("N41").Value = "$/caja 16 unidades"
$ws.Range("O41").Value = "Provincia de Quillota"
$ws.Range("P41").Value = 384
$ws.Range("Q41").Value = 16
$ws.Range("R41").Value = "Hortaliza"
